$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new user record (row 3) below the existing username/password rows
$ws.Range("A3").Value = "Ruby"
$ws.Range("B3").Value = "22333LJ"

# Match the author's final selection, sitting on the newly added cell
$ws.Range("B3").Select()
